$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, C->D)
$ws.Columns("B:B").Insert()

# Match the width of the new column B to column A's width (both are
# 75.81640625 "characters" wide in the target file; the closest value this
# host's pixel-quantized ColumnWidth model can reproduce is 75).
$ws.Columns("B").ColumnWidth = 75

# New header + query cells
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE a.pubmed_id IN ['31504139'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# B2 should wrap text like A2 (it already inherits this from the Insert,
# but set explicitly to be certain)
$ws.Range("B2").WrapText = $true

# Move the active selection to B4, matching the saved view state
$ws.Range("B4").Select() | Out-Null
